$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C9").Value = "Exponencial"
$excel.ActiveWindow.Zoom = 280
$ws.Range("B1:D9").Select() | Out-Null
$ws.Rows.Item(9).AutoFit() | Out-Null
